$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''27.035.78'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '''  -1.03%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = '''1.822.98'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '''  -0.62%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').Value = '''1.007'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '''  -0.36%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = '''309.95'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '''  -1.53%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('E6').Value = '''  -0.27%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = '''0.4652'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '''  -1.88%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = '''0.3664'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '''  -0.72%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = '''0.07243'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '''  -2.71%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = '''0.8595'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '''  -2.98%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('B11').Value = '''WrappedEther'
$ws.Range('B11').Style = 'Normal'
$ws.Range('C11').Value = '''https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('C11').Style = 'Normal'
$ws.Range('D11').Value = '''1.933.96'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '''  +3.49%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('B12').Value = '''Solana'
$ws.Range('B12').Style = 'Normal'
$ws.Range('C12').Value = '''https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('C12').Style = 'Normal'
$ws.Range('D12').Value = '''19.86'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '''  -2.96%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('B13').Value = '''TRON'
$ws.Range('B13').Style = 'Normal'
$ws.Range('C13').Value = '''https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('C13').Style = 'Normal'
$ws.Range('D13').Value = '''0.07704'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '''  +5.10%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = '''5.329'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '''  -2.05%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('B15').Value = '''Litecoin'
$ws.Range('B15').Style = 'Normal'
$ws.Range('C15').Value = '''https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('C15').Style = 'Normal'
$ws.Range('D15').Value = '''91.75'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '''  -2.26%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('B16').Value = '''Chainlink'
$ws.Range('B16').Style = 'Normal'
$ws.Range('C16').Value = '''https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('C16').Style = 'Normal'
$ws.Range('D16').Value = '''6.499'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '''  -1.17%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('E17').Value = '''  -0.24%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = '''0.000008655'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '''  -1.51%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = '''1.007'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '''  -0.37%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = '''26.938.45'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '''  -2.08%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = '''14.48'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Value = '''5.147'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '''  -2.69%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = '''10.54'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '''  -1.20%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = '''2.086.17'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '''  -0.10%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = '''151.85'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '''  -0.12%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = '''1.840'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '''  -2.91%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('D28').Value = '''2.047'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '''  -4.84%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = '''5.097'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '''  -2.54%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = '''115.38'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '''  -1.58%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = '''0.08841'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '''  -1.69%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('E32').Value = '''  +0.22%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = '''4.429'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '''  -2.61%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = '''1.130'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '''  -3.90%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = '''0.7206'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '''  -4.07%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = '''1.076'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '''  -2.28%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = '''0.05248'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '''  -1.81%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = '''0.01925'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '''  -1.54%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = '''2.406'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '''  +0.80%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = '''2.930'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '''  -1.33%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = '''7.122'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '''  -1.49%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = '''0.5164'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '''  -2.79%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('B43').Value = '''Frax'
$ws.Range('B43').Style = 'Normal'
$ws.Range('C43').Value = '''https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('C43').Style = 'Normal'
$ws.Range('D43').Value = '''0.8720'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '''  -13.68%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('B44').Value = '''Algorand'
$ws.Range('B44').Style = 'Normal'
$ws.Range('C44').Value = '''https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('C44').Style = 'Normal'
$ws.Range('D44').Value = '''0.1628'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '''  -1.94%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = '''8.166'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '''  -3.77%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = '''0.4799'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '''  -2.77%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('E47').Value = '''  -0.27%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = '''10.10'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '''  -4.41%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = '''102.57'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '''  -2.51%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = '''0.06249'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '''  -0.75%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = '''1.618'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '''  -3.31%  '
$ws.Range('E51').Style = 'Normal'
